$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 'Achu'
$ws.Range("B3").Value = 'ashwinor0812@gmail.com'
$ws.Range("C3").Value = 'scrypt:32768:8:1$FdWQWR7qY0iQLMXR$bde4fd2a777f18096b4a1575d8bd892f2dce9041e4ed6dd3f7573c68f6fb26adc8274d46a0312448f77573885bd5df6073d4425328c177b6d0bd7e5a774561f9'

# Row 4
$ws.Range("A4").Value = 'testuser2'
$ws.Range("B4").Value = 'test2@test.com'
$ws.Range("C4").Value = 'testpass'

# Row 5
$ws.Range("A5").Value = 'ashw'
$ws.Range("B5").Value = 'ashwinor000@gmail.com'
$ws.Range("C5").Value = 'scrypt:32768:8:1$ffvEDCcw7JImFPBH$c13c0b0eaead011e78a41f3fe2e2adca7b95f99a9ab2e4e1647942209d5dd2e009e7fcdb28d509f939de53ea09812cc33c550859bf5386dd002ac93aef5c6717'
